# "corrected multiple worksheets articles"
# The workbook is a document-generation demo template. The merge-field
# placeholder used for the worksheet/tab name was renamed from `{{name}}`
# to the more specific `{{sheetName}}`, and the sheet's saved selection
# was reset to the header row (A1:B1) instead of a stray mid-sheet cell.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rename the worksheet tab - this is the `{{name}}` placeholder that shows
# up as the <sheet name="..."/> entry in xl/workbook.xml.
$ws.Name = "{{sheetName}}"

# Reset the sheet's active selection to the header row range A1:B1
# (previously parked on C9).
$ws.Range("A1:B1").Select()
